$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
if ($ws.Range("F3").Value2 -ne 238) { throw "Unexpected value in 展览!F3: expected 238, got $($ws.Range('F3').Value2)" }
$ws.Range("F3").Value2 = 240
if ($ws.Range("F5").Value2 -ne 25) { throw "Unexpected value in 展览!F5: expected 25, got $($ws.Range('F5').Value2)" }
$ws.Range("F5").Value2 = 27
if ($ws.Range("F6").Value2 -ne 93) { throw "Unexpected value in 展览!F6: expected 93, got $($ws.Range('F6').Value2)" }
$ws.Range("F6").Value2 = 103
if ($ws.Range("F7").Value2 -ne 12) { throw "Unexpected value in 展览!F7: expected 12, got $($ws.Range('F7').Value2)" }
$ws.Range("F7").Value2 = 13
if ($ws.Range("F9").Value2 -ne 4648) { throw "Unexpected value in 展览!F9: expected 4648, got $($ws.Range('F9').Value2)" }
$ws.Range("F9").Value2 = 4668
if ($ws.Range("F10").Value2 -ne 4648) { throw "Unexpected value in 展览!F10: expected 4648, got $($ws.Range('F10').Value2)" }
$ws.Range("F10").Value2 = 4668
if ($ws.Range("F12").Value2 -ne 446) { throw "Unexpected value in 展览!F12: expected 446, got $($ws.Range('F12').Value2)" }
$ws.Range("F12").Value2 = 447
if ($ws.Range("F13").Value2 -ne 1084) { throw "Unexpected value in 展览!F13: expected 1084, got $($ws.Range('F13').Value2)" }
$ws.Range("F13").Value2 = 1088
if ($ws.Range("F14").Value2 -ne 605) { throw "Unexpected value in 展览!F14: expected 605, got $($ws.Range('F14').Value2)" }
$ws.Range("F14").Value2 = 608
if ($ws.Range("F15").Value2 -ne 4149) { throw "Unexpected value in 展览!F15: expected 4149, got $($ws.Range('F15').Value2)" }
$ws.Range("F15").Value2 = 4186
if ($ws.Range("F16").Value2 -ne 160) { throw "Unexpected value in 展览!F16: expected 160, got $($ws.Range('F16').Value2)" }
$ws.Range("F16").Value2 = 161
if ($ws.Range("F17").Value2 -ne 161) { throw "Unexpected value in 展览!F17: expected 161, got $($ws.Range('F17').Value2)" }
$ws.Range("F17").Value2 = 162
if ($ws.Range("F18").Value2 -ne 48) { throw "Unexpected value in 展览!F18: expected 48, got $($ws.Range('F18').Value2)" }
$ws.Range("F18").Value2 = 49
if ($ws.Range("F19").Value2 -ne 203) { throw "Unexpected value in 展览!F19: expected 203, got $($ws.Range('F19').Value2)" }
$ws.Range("F19").Value2 = 209
if ($ws.Range("F20").Value2 -ne 3424) { throw "Unexpected value in 展览!F20: expected 3424, got $($ws.Range('F20').Value2)" }
$ws.Range("F20").Value2 = 3439
if ($ws.Range("F24").Value2 -ne 2988) { throw "Unexpected value in 展览!F24: expected 2988, got $($ws.Range('F24').Value2)" }
$ws.Range("F24").Value2 = 3021
if ($ws.Range("F25").Value2 -ne 124) { throw "Unexpected value in 展览!F25: expected 124, got $($ws.Range('F25').Value2)" }
$ws.Range("F25").Value2 = 125
if ($ws.Range("F26").Value2 -ne 124) { throw "Unexpected value in 展览!F26: expected 124, got $($ws.Range('F26').Value2)" }
$ws.Range("F26").Value2 = 125
if ($ws.Range("F28").Value2 -ne 143) { throw "Unexpected value in 展览!F28: expected 143, got $($ws.Range('F28').Value2)" }
$ws.Range("F28").Value2 = 145
if ($ws.Range("F29").Value2 -ne 180) { throw "Unexpected value in 展览!F29: expected 180, got $($ws.Range('F29').Value2)" }
$ws.Range("F29").Value2 = 184
if ($ws.Range("F30").Value2 -ne 172) { throw "Unexpected value in 展览!F30: expected 172, got $($ws.Range('F30').Value2)" }
$ws.Range("F30").Value2 = 174
if ($ws.Range("F31").Value2 -ne 72) { throw "Unexpected value in 展览!F31: expected 72, got $($ws.Range('F31').Value2)" }
$ws.Range("F31").Value2 = 74
if ($ws.Range("F32").Value2 -ne 49) { throw "Unexpected value in 展览!F32: expected 49, got $($ws.Range('F32').Value2)" }
$ws.Range("F32").Value2 = 51
if ($ws.Range("F33").Value2 -ne 22) { throw "Unexpected value in 展览!F33: expected 22, got $($ws.Range('F33').Value2)" }
$ws.Range("F33").Value2 = 23
if ($ws.Range("F36").Value2 -ne 5380) { throw "Unexpected value in 展览!F36: expected 5380, got $($ws.Range('F36').Value2)" }
$ws.Range("F36").Value2 = 5420
if ($ws.Range("F37").Value2 -ne 752) { throw "Unexpected value in 展览!F37: expected 752, got $($ws.Range('F37').Value2)" }
$ws.Range("F37").Value2 = 767
if ($ws.Range("F38").Value2 -ne 384) { throw "Unexpected value in 展览!F38: expected 384, got $($ws.Range('F38').Value2)" }
$ws.Range("F38").Value2 = 386
if ($ws.Range("F39").Value2 -ne 80) { throw "Unexpected value in 展览!F39: expected 80, got $($ws.Range('F39').Value2)" }
$ws.Range("F39").Value2 = 82
if ($ws.Range("F41").Value2 -ne 31) { throw "Unexpected value in 展览!F41: expected 31, got $($ws.Range('F41').Value2)" }
$ws.Range("F41").Value2 = 43
if ($ws.Range("F42").Value2 -ne 1072) { throw "Unexpected value in 展览!F42: expected 1072, got $($ws.Range('F42').Value2)" }
$ws.Range("F42").Value2 = 1085
if ($ws.Range("F43").Value2 -ne 447) { throw "Unexpected value in 展览!F43: expected 447, got $($ws.Range('F43').Value2)" }
$ws.Range("F43").Value2 = 456
if ($ws.Range("F45").Value2 -ne 1942) { throw "Unexpected value in 展览!F45: expected 1942, got $($ws.Range('F45').Value2)" }
$ws.Range("F45").Value2 = 1951
if ($ws.Range("F47").Value2 -ne 53) { throw "Unexpected value in 展览!F47: expected 53, got $($ws.Range('F47').Value2)" }
$ws.Range("F47").Value2 = 56
if ($ws.Range("F48").Value2 -ne 691) { throw "Unexpected value in 展览!F48: expected 691, got $($ws.Range('F48').Value2)" }
$ws.Range("F48").Value2 = 692
if ($ws.Range("F49").Value2 -ne 828) { throw "Unexpected value in 展览!F49: expected 828, got $($ws.Range('F49').Value2)" }
$ws.Range("F49").Value2 = 835

$ws = $wb.Worksheets.Item("演出")
if ($ws.Range("F3").Value2 -ne 8) { throw "Unexpected value in 演出!F3: expected 8, got $($ws.Range('F3').Value2)" }
$ws.Range("F3").Value2 = 9
if ($ws.Range("F8").Value2 -ne 37) { throw "Unexpected value in 演出!F8: expected 37, got $($ws.Range('F8').Value2)" }
$ws.Range("F8").Value2 = 38
if ($ws.Range("F15").Value2 -ne 119) { throw "Unexpected value in 演出!F15: expected 119, got $($ws.Range('F15').Value2)" }
$ws.Range("F15").Value2 = 120
if ($ws.Range("F22").Value2 -ne 722) { throw "Unexpected value in 演出!F22: expected 722, got $($ws.Range('F22').Value2)" }
$ws.Range("F22").Value2 = 727

$ws = $wb.Worksheets.Item("全部类型")
if ($ws.Range("F3").Value2 -ne 8) { throw "Unexpected value in 全部类型!F3: expected 8, got $($ws.Range('F3').Value2)" }
$ws.Range("F3").Value2 = 9
if ($ws.Range("F5").Value2 -ne 238) { throw "Unexpected value in 全部类型!F5: expected 238, got $($ws.Range('F5').Value2)" }
$ws.Range("F5").Value2 = 240
if ($ws.Range("F6").Value2 -ne 25) { throw "Unexpected value in 全部类型!F6: expected 25, got $($ws.Range('F6').Value2)" }
$ws.Range("F6").Value2 = 27
if ($ws.Range("F8").Value2 -ne 93) { throw "Unexpected value in 全部类型!F8: expected 93, got $($ws.Range('F8').Value2)" }
$ws.Range("F8").Value2 = 103
if ($ws.Range("F9").Value2 -ne 12) { throw "Unexpected value in 全部类型!F9: expected 12, got $($ws.Range('F9').Value2)" }
$ws.Range("F9").Value2 = 13
if ($ws.Range("F11").Value2 -ne 4648) { throw "Unexpected value in 全部类型!F11: expected 4648, got $($ws.Range('F11').Value2)" }
$ws.Range("F11").Value2 = 4668
if ($ws.Range("F12").Value2 -ne 4648) { throw "Unexpected value in 全部类型!F12: expected 4648, got $($ws.Range('F12').Value2)" }
$ws.Range("F12").Value2 = 4668
if ($ws.Range("F13").Value2 -ne 37) { throw "Unexpected value in 全部类型!F13: expected 37, got $($ws.Range('F13').Value2)" }
$ws.Range("F13").Value2 = 38
if ($ws.Range("F17").Value2 -ne 446) { throw "Unexpected value in 全部类型!F17: expected 446, got $($ws.Range('F17').Value2)" }
$ws.Range("F17").Value2 = 447
if ($ws.Range("F18").Value2 -ne 1084) { throw "Unexpected value in 全部类型!F18: expected 1084, got $($ws.Range('F18').Value2)" }
$ws.Range("F18").Value2 = 1088
if ($ws.Range("F19").Value2 -ne 605) { throw "Unexpected value in 全部类型!F19: expected 605, got $($ws.Range('F19').Value2)" }
$ws.Range("F19").Value2 = 608
if ($ws.Range("F20").Value2 -ne 4149) { throw "Unexpected value in 全部类型!F20: expected 4149, got $($ws.Range('F20').Value2)" }
$ws.Range("F20").Value2 = 4186
if ($ws.Range("F21").Value2 -ne 160) { throw "Unexpected value in 全部类型!F21: expected 160, got $($ws.Range('F21').Value2)" }
$ws.Range("F21").Value2 = 161
if ($ws.Range("F22").Value2 -ne 161) { throw "Unexpected value in 全部类型!F22: expected 161, got $($ws.Range('F22').Value2)" }
$ws.Range("F22").Value2 = 162
if ($ws.Range("F23").Value2 -ne 203) { throw "Unexpected value in 全部类型!F23: expected 203, got $($ws.Range('F23').Value2)" }
$ws.Range("F23").Value2 = 209
if ($ws.Range("F24").Value2 -ne 3424) { throw "Unexpected value in 全部类型!F24: expected 3424, got $($ws.Range('F24').Value2)" }
$ws.Range("F24").Value2 = 3439
if ($ws.Range("F25").Value2 -ne 2988) { throw "Unexpected value in 全部类型!F25: expected 2988, got $($ws.Range('F25').Value2)" }
$ws.Range("F25").Value2 = 3021
if ($ws.Range("F26").Value2 -ne 124) { throw "Unexpected value in 全部类型!F26: expected 124, got $($ws.Range('F26').Value2)" }
$ws.Range("F26").Value2 = 125
if ($ws.Range("F27").Value2 -ne 124) { throw "Unexpected value in 全部类型!F27: expected 124, got $($ws.Range('F27').Value2)" }
$ws.Range("F27").Value2 = 125
if ($ws.Range("F28").Value2 -ne 143) { throw "Unexpected value in 全部类型!F28: expected 143, got $($ws.Range('F28').Value2)" }
$ws.Range("F28").Value2 = 145
if ($ws.Range("F29").Value2 -ne 180) { throw "Unexpected value in 全部类型!F29: expected 180, got $($ws.Range('F29').Value2)" }
$ws.Range("F29").Value2 = 184
if ($ws.Range("F30").Value2 -ne 172) { throw "Unexpected value in 全部类型!F30: expected 172, got $($ws.Range('F30').Value2)" }
$ws.Range("F30").Value2 = 174
if ($ws.Range("F31").Value2 -ne 22) { throw "Unexpected value in 全部类型!F31: expected 22, got $($ws.Range('F31').Value2)" }
$ws.Range("F31").Value2 = 23
if ($ws.Range("F35").Value2 -ne 119) { throw "Unexpected value in 全部类型!F35: expected 119, got $($ws.Range('F35').Value2)" }
$ws.Range("F35").Value2 = 120
if ($ws.Range("F37").Value2 -ne 5380) { throw "Unexpected value in 全部类型!F37: expected 5380, got $($ws.Range('F37').Value2)" }
$ws.Range("F37").Value2 = 5420
if ($ws.Range("F39").Value2 -ne 752) { throw "Unexpected value in 全部类型!F39: expected 752, got $($ws.Range('F39').Value2)" }
$ws.Range("F39").Value2 = 767
if ($ws.Range("F40").Value2 -ne 384) { throw "Unexpected value in 全部类型!F40: expected 384, got $($ws.Range('F40').Value2)" }
$ws.Range("F40").Value2 = 386
if ($ws.Range("F42").Value2 -ne 80) { throw "Unexpected value in 全部类型!F42: expected 80, got $($ws.Range('F42').Value2)" }
$ws.Range("F42").Value2 = 82
if ($ws.Range("F44").Value2 -ne 1072) { throw "Unexpected value in 全部类型!F44: expected 1072, got $($ws.Range('F44').Value2)" }
$ws.Range("F44").Value2 = 1085
if ($ws.Range("F45").Value2 -ne 447) { throw "Unexpected value in 全部类型!F45: expected 447, got $($ws.Range('F45').Value2)" }
$ws.Range("F45").Value2 = 456
if ($ws.Range("F47").Value2 -ne 1942) { throw "Unexpected value in 全部类型!F47: expected 1942, got $($ws.Range('F47').Value2)" }
$ws.Range("F47").Value2 = 1951
if ($ws.Range("F48").Value2 -ne 53) { throw "Unexpected value in 全部类型!F48: expected 53, got $($ws.Range('F48').Value2)" }
$ws.Range("F48").Value2 = 56
if ($ws.Range("F49").Value2 -ne 691) { throw "Unexpected value in 全部类型!F49: expected 691, got $($ws.Range('F49').Value2)" }
$ws.Range("F49").Value2 = 692
if ($ws.Range("F50").Value2 -ne 828) { throw "Unexpected value in 全部类型!F50: expected 828, got $($ws.Range('F50').Value2)" }
$ws.Range("F50").Value2 = 835
